$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos
$ws.Range("C3").Value = "El Rubio"
$ws.Range("B25").Value = "Amigo de Kibelo"

# Add new column D: avatar_img
$ws.Range("D1").Value = "avatar_img"
$ws.Range("D2").Value = "johan.jpg"
$ws.Range("D6").Value = "randy.jpg"
$ws.Range("D7").Value = "punto.jpg"
$ws.Range("D9").Value = "joel.jpg"
$ws.Range("D10").Value = "cristopher.jpg"
$ws.Range("D13").Value = "kukito.jpg"
$ws.Range("D15").Value = "alfredo.jpg"
$ws.Range("D21").Value = "omauri.jpg"
$ws.Range("D27").Value = "carlos.jpg"
$ws.Range("D28").Value = "kawai.jpg"
$ws.Range("D30").Value = "rayder.jpg"
$ws.Range("D32").Value = "jeicol.jpg"
$ws.Range("D34").Value = "yeyo.jpg"
